# Atualização automática SALDO_PECAS (18/11/2025 12:43)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C12: was stored as the text "125200", becomes a real number 125200 ---
$ws.Range("C12").Value = 125200
$ws.Range("C12").Style = "Normal"

# --- New row 13 ---
# Write every value as literal text via the classic leading single-quote
# ("force text") prefix so strings that look numeric/date-like (e.g.
# "0943000", "01/11/25") are never silently reinterpreted as a number or
# date, then drop the cell back to the default "Normal" style so no
# quote-prefix / number-format marker is left behind on the cell.
function Set-TextCell($addr, $value) {
    $ws.Range($addr).Value = "'" + $value
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell "A13" "DF"
Set-TextCell "B13" "DF18110"
Set-TextCell "C13" "0943000"
Set-TextCell "D13" ""
Set-TextCell "E13" ""
Set-TextCell "F13" "T"
Set-TextCell "G13" "T"
Set-TextCell "H13" "T - (T 01/11/25_24H) - DF"
Set-TextCell "I13" "01/11/25"
Set-TextCell "J13" "24H"
Set-TextCell "K13" "18/11/25"
Set-TextCell "L13" "DENTRO"
Set-TextCell "M13" ""
